# Agrega proceso tercierizado de cuño
# Inserts a new "Cuño" (tercerizado) process row into the "Maquinas" sheet,
# pushing the existing Troquelado/Descartonado/Ventana/Pegado rows down by one,
# and updates the active-sheet / selection state to match.

$wb = $excel.ActiveWorkbook

$wsMaquinas = $wb.Worksheets.Item("Maquinas")
$wsOrden = $wb.Worksheets.Item("OrdenEstandar")

# Insert a new row at row 9 (shifts Troquelado/Descartonado/Ventana/Pegado rows down)
$wsMaquinas.Rows.Item(9).Insert()

# Populate the new "Cuño" row (Proceso = Cuño, Maquina = Cuño)
$wsMaquinas.Cells.Item(9, 1).Value = "Cuño"
$wsMaquinas.Cells.Item(9, 2).Value = "Cuño"
$wsMaquinas.Cells.Item(9, 3).Value = 1000
$wsMaquinas.Cells.Item(9, 4).Value = 10
$wsMaquinas.Cells.Item(9, 5).Value = 10

# Update selection on OrdenEstandar (no longer the active tab)
$wsOrden.Range("B10").Select()

# Make Maquinas the active sheet with the new selection
$wsMaquinas.Activate()
$wsMaquinas.Range("C9").Select()
